# ecYaliGEM/data/customKcats.xlsx update
# "New ecYali based on an update of iYali4"
#
# Adds three new custom-kcat entries (rows 43-45) to the "customKcats"
# sheet, and updates the active selection to reflect where the author
# ended up after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customKcats")
$ws.Activate()

# --- Row 43: YALI0B02948g / Q9UVF4 --------------------------------------
$ws.Range("A43").Value = "Q9UVF4"
$ws.Range("B43").Value = "YALI0B02948g"
$ws.Range("C43").Value = "YALI0B02948g"
$ws.Range("D43").Value = 114.442
$ws.Range("E43").Value = "y000491_REV"
$ws.Range("F43").Value = "Limits model after curations of iYali. Calculated from the specific activity of S. cerevisiae (EC 1.1.1.8)"
$ws.Range("G43").Value = 1

# --- Row 44: YALI0D16753g / Q6C8V3 (note added further below) ----------
$ws.Range("A44").Value = "Q6C8V3"
$ws.Range("B44").Value = "YALI0D16753g"
$ws.Range("C44").Value = "YALI0D16753g"
$ws.Range("D44").Value = 83.2
$ws.Range("E44").Value = "y000713"
$ws.Range("G44").Value = 1

# --- Row 45: YALI0E14190g / Q6C5X9 --------------------------------------
$ws.Range("A45").Value = "Q6C5X9"
$ws.Range("B45").Value = "YALI0E14190g"
$ws.Range("C45").Value = "YALI0E14190g"
$ws.Range("D45").Value = 3.1018
$ws.Range("E45").Value = "y000713"
$ws.Range("F45").Value = "kcat from BRENDA was causing overuse of the NADH shuttle. Using value of DLKcat instead."
$ws.Range("G45").Value = 1

# Note for row 44 filled in last, matching the authoring order recorded
# in the shared-strings table.
$ws.Range("F44").Value = "kcat from BRENDA (highest) was causing overuse of the NADH shuttle. Using highest kcat in the same order of magnitude predicted by DLKcat."

# Leave the selection where the author left it.
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("A37").Select()
